$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 3053.5405
$ws.Cells.Item(15, 9).Value = 3053.5405
$ws.Cells.Item(15, 11).Value = 9160.621500000001
$ws.Cells.Item(15, 13).Value = -8991.621500000001
$ws.Cells.Item(33, 8).Value = 193.17647
$ws.Cells.Item(33, 9).Value = 140.4
$ws.Cells.Item(33, 10).Value = 589
$ws.Cells.Item(33, 11).Value = 140.4
$ws.Cells.Item(33, 12).Value = 589
$ws.Cells.Item(33, 13).Value = 88.59999999999999
$ws.Cells.Item(33, 14).Value = -1047
$ws.Cells.Item(135, 8).Value = 40000644
$ws.Cells.Item(135, 9).Value = 321.38095
$ws.Cells.Item(135, 10).Value = 250002340
$ws.Cells.Item(135, 11).Value = 2892.42855
$ws.Cells.Item(135, 12).Value = 2250021060
$ws.Cells.Item(135, 13).Value = -357.4285499999996
$ws.Cells.Item(135, 14).Value = -2250026130
$ws.Cells.Item(137, 8).Value = 981.50616
$ws.Cells.Item(137, 9).Value = 664.92
$ws.Cells.Item(137, 11).Value = 1994.76
$ws.Cells.Item(137, 13).Value = 555.2400000000002
$ws.Cells.Item(138, 8).Value = 1160.51
$ws.Cells.Item(138, 9).Value = 524.7954999999999
$ws.Cells.Item(138, 10).Value = 1660
$ws.Cells.Item(138, 11).Value = 1574.3865
$ws.Cells.Item(138, 12).Value = 4980
$ws.Cells.Item(138, 13).Value = 3565.6135
$ws.Cells.Item(138, 14).Value = -15260
$ws.Cells.Item(141, 8).Value = 419.65
$ws.Cells.Item(141, 9).Value = 441
$ws.Cells.Item(141, 10).Value = 227.5
$ws.Cells.Item(141, 11).Value = 1323
$ws.Cells.Item(141, 12).Value = 682.5
$ws.Cells.Item(141, 13).Value = 3857
$ws.Cells.Item(141, 14).Value = -11042.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5466.1924
$ws.Cells.Item(32, 9).Value = 4631.6562
$ws.Cells.Item(32, 10).Value = 9281.214
$ws.Cells.Item(32, 11).Value = 4631.6562
$ws.Cells.Item(32, 12).Value = 9281.214
$ws.Cells.Item(32, 13).Value = -4344.6562
$ws.Cells.Item(32, 14).Value = -9855.214
$ws.Cells.Item(61, 8).Value = 47620060
$ws.Cells.Item(61, 9).Value = 58824444
$ws.Cells.Item(61, 10).Value = 1425
$ws.Cells.Item(61, 11).Value = 58824444
$ws.Cells.Item(61, 12).Value = 1425
$ws.Cells.Item(61, 13).Value = -58824232
$ws.Cells.Item(61, 14).Value = -1849
$ws.Cells.Item(88, 8).Value = 2539.6
$ws.Cells.Item(88, 10).Value = 2716.1667
$ws.Cells.Item(88, 12).Value = 2716.1667
$ws.Cells.Item(88, 14).Value = -3528.1667
$ws.Cells.Item(91, 8).Value = 2539.6
$ws.Cells.Item(91, 10).Value = 2716.1667
$ws.Cells.Item(91, 12).Value = 2716.1667
$ws.Cells.Item(91, 14).Value = -5524.1667
$ws.Cells.Item(97, 8).Value = 291.84848
$ws.Cells.Item(97, 9).Value = 308.80768
$ws.Cells.Item(97, 10).Value = 228.85715
$ws.Cells.Item(97, 11).Value = 308.80768
$ws.Cells.Item(97, 12).Value = 228.85715
$ws.Cells.Item(97, 13).Value = 187.19232
$ws.Cells.Item(97, 14).Value = -1220.85715
$ws.Cells.Item(136, 8).Value = 47620060
$ws.Cells.Item(136, 9).Value = 58824444
$ws.Cells.Item(136, 10).Value = 1425
$ws.Cells.Item(136, 11).Value = 176473332
$ws.Cells.Item(136, 12).Value = 4275
$ws.Cells.Item(136, 13).Value = -176470782
$ws.Cells.Item(136, 14).Value = -9375
$ws.Cells.Item(141, 8).Value = 32066.334
$ws.Cells.Item(141, 10).Value = 32066.334
$ws.Cells.Item(141, 12).Value = 32066.334
$ws.Cells.Item(141, 14).Value = -42426.334

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 856.4186
$ws.Cells.Item(107, 9).Value = 689.6111
$ws.Cells.Item(107, 11).Value = 689.6111
$ws.Cells.Item(107, 13).Value = 1230.3889
$ws.Cells.Item(134, 8).Value = 2920.8245
$ws.Cells.Item(134, 9).Value = 1014.55554
$ws.Cells.Item(134, 10).Value = 10069.333
$ws.Cells.Item(134, 11).Value = 3043.66662
$ws.Cells.Item(134, 12).Value = 30207.999
$ws.Cells.Item(134, 13).Value = -508.66662
$ws.Cells.Item(134, 14).Value = -35277.999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1840.6
$ws.Cells.Item(31, 9).Value = 1896.84
$ws.Cells.Item(31, 10).Value = 1559.4
$ws.Cells.Item(31, 11).Value = 1896.84
$ws.Cells.Item(31, 12).Value = 1559.4
$ws.Cells.Item(31, 13).Value = -1601.84
$ws.Cells.Item(31, 14).Value = -2149.4
$ws.Cells.Item(34, 8).Value = 1840.6
$ws.Cells.Item(34, 9).Value = 1896.84
$ws.Cells.Item(34, 10).Value = 1559.4
$ws.Cells.Item(34, 11).Value = 1896.84
$ws.Cells.Item(34, 12).Value = 1559.4
$ws.Cells.Item(34, 13).Value = -1694.84
$ws.Cells.Item(34, 14).Value = -1963.4
$ws.Cells.Item(58, 8).Value = 917.56604
$ws.Cells.Item(58, 9).Value = 872.75
$ws.Cells.Item(58, 10).Value = 1136.6666
$ws.Cells.Item(58, 11).Value = 872.75
$ws.Cells.Item(58, 12).Value = 1136.6666
$ws.Cells.Item(58, 13).Value = -669.75
$ws.Cells.Item(58, 14).Value = -1542.6666
$ws.Cells.Item(132, 8).Value = 3327.5
$ws.Cells.Item(132, 9).Value = 4086.5
$ws.Cells.Item(132, 10).Value = 1809.5
$ws.Cells.Item(132, 11).Value = 12259.5
$ws.Cells.Item(132, 12).Value = 5428.5
$ws.Cells.Item(132, 13).Value = -9729.5
$ws.Cells.Item(132, 14).Value = -10488.5
$ws.Cells.Item(134, 8).Value = 9435132
$ws.Cells.Item(134, 9).Value = 1209
$ws.Cells.Item(134, 10).Value = 31251078
$ws.Cells.Item(134, 11).Value = 3627
$ws.Cells.Item(134, 12).Value = 93753234
$ws.Cells.Item(134, 13).Value = -1092
$ws.Cells.Item(134, 14).Value = -93758304
$ws.Cells.Item(136, 8).Value = 917.56604
$ws.Cells.Item(136, 9).Value = 872.75
$ws.Cells.Item(136, 10).Value = 1136.6666
$ws.Cells.Item(136, 11).Value = 2618.25
$ws.Cells.Item(136, 12).Value = 3409.9998
$ws.Cells.Item(136, 13).Value = -68.25
$ws.Cells.Item(136, 14).Value = -8509.9998
$ws.Cells.Item(141, 8).Value = 76163.71000000001
$ws.Cells.Item(141, 10).Value = 76163.71000000001
$ws.Cells.Item(141, 12).Value = 76163.71000000001
$ws.Cells.Item(141, 14).Value = -86523.71000000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 177.91304
$ws.Cells.Item(14, 9).Value = 177.91304
$ws.Cells.Item(14, 11).Value = 533.73912
$ws.Cells.Item(14, 13).Value = -360.73912
$ws.Cells.Item(131, 8).Value = 20836358
$ws.Cells.Item(131, 9).Value = 76923560
$ws.Cells.Item(131, 10).Value = 3967.4
$ws.Cells.Item(131, 11).Value = 230770680
$ws.Cells.Item(131, 12).Value = 11902.2
$ws.Cells.Item(131, 13).Value = -230765640
$ws.Cells.Item(131, 14).Value = -21982.2
$ws.Cells.Item(139, 8).Value = 2044.4103
$ws.Cells.Item(139, 9).Value = 2333.7273
$ws.Cells.Item(139, 10).Value = 1670
$ws.Cells.Item(139, 11).Value = 7001.1819
$ws.Cells.Item(139, 12).Value = 5010
$ws.Cells.Item(139, 13).Value = -1861.1819
$ws.Cells.Item(139, 14).Value = -15290

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 565.17645
$ws.Cells.Item(102, 9).Value = 579.1429000000001
$ws.Cells.Item(102, 11).Value = 579.1429000000001
$ws.Cells.Item(102, 13).Value = 1042.8571

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 490.125
$ws.Cells.Item(22, 9).Value = 434.5
$ws.Cells.Item(22, 11).Value = 434.5
$ws.Cells.Item(22, 13).Value = -139.5
$ws.Cells.Item(27, 8).Value = 490.125
$ws.Cells.Item(27, 9).Value = 434.5
$ws.Cells.Item(27, 11).Value = 434.5
$ws.Cells.Item(27, 13).Value = -327.5
$ws.Cells.Item(68, 8).Value = 1913.375
$ws.Cells.Item(68, 9).Value = 1916.9474
$ws.Cells.Item(68, 10).Value = 1899.8
$ws.Cells.Item(68, 11).Value = 1916.9474
$ws.Cells.Item(68, 12).Value = 1899.8
$ws.Cells.Item(68, 13).Value = -1167.9474
$ws.Cells.Item(68, 14).Value = -3397.8
$ws.Cells.Item(71, 8).Value = 1913.375
$ws.Cells.Item(71, 9).Value = 1916.9474
$ws.Cells.Item(71, 10).Value = 1899.8
$ws.Cells.Item(71, 11).Value = 9584.737000000001
$ws.Cells.Item(71, 12).Value = 9499
$ws.Cells.Item(71, 13).Value = -5840.737000000001
$ws.Cells.Item(71, 14).Value = -16987
$ws.Cells.Item(132, 8).Value = 24694.705
$ws.Cells.Item(132, 9).Value = 1502.4073
$ws.Cells.Item(132, 10).Value = 61529.53
$ws.Cells.Item(132, 11).Value = 4507.2219
$ws.Cells.Item(132, 12).Value = 184588.59
$ws.Cells.Item(132, 13).Value = -1977.2219
$ws.Cells.Item(132, 14).Value = -189648.59
$ws.Cells.Item(136, 8).Value = 2847.6345
$ws.Cells.Item(136, 9).Value = 2938.6956
$ws.Cells.Item(136, 11).Value = 8816.086800000001
$ws.Cells.Item(136, 13).Value = -6266.086800000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 65361216
$ws.Cells.Item(126, 9).Value = 69446104
$ws.Cells.Item(126, 11).Value = 208338312
$ws.Cells.Item(126, 13).Value = -208335842
$ws.Cells.Item(132, 8).Value = 1798.9736
$ws.Cells.Item(132, 9).Value = 1609.9822
$ws.Cells.Item(132, 10).Value = 2328.15
$ws.Cells.Item(132, 11).Value = 4829.946599999999
$ws.Cells.Item(132, 12).Value = 6984.450000000001
$ws.Cells.Item(132, 13).Value = -2299.946599999999
$ws.Cells.Item(132, 14).Value = -12044.45
$ws.Cells.Item(136, 8).Value = 569.3182
$ws.Cells.Item(136, 9).Value = 509.0909
$ws.Cells.Item(136, 10).Value = 750
$ws.Cells.Item(136, 11).Value = 1527.2727
$ws.Cells.Item(136, 12).Value = 2250
$ws.Cells.Item(136, 13).Value = 1022.7273
$ws.Cells.Item(136, 14).Value = -7350
